$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - Planung: mark 100% complete
$ws.Range("G5").Value = 1

# Row 8 - Planungs-Doku: mark 100% complete
$ws.Range("G8").Value = 1

# Row 9 - Implementierung: fill in actual start/duration, mark complete
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 7
$ws.Range("G9").Value = 1

# Row 11 - Szenen-Objekte
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1

# Row 12 - Kamera-Steuerung
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 1

# Row 13 - Update-Gespräch
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1

# Row 14 - Dokumentation
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 1

# Row 15 - Zeit-Doku
$ws.Range("E15").Value = 11
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1

# Row 16 - Präsentations-Erstellung
$ws.Range("E16").Value = 11
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1

# Row 17 - Andwender-Doku
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1

# Row 18 - Projekt-Doku
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 1

# Row 19 - Abgabe
$ws.Range("E19").Value = 11
$ws.Range("F19").Value = 2
$ws.Range("G19").Value = 1

# Row 20 - Finaler Upload
$ws.Range("E20").Value = 11
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1

# Row 21 - Präsentation (F21 stays as "?" placeholder, unchanged)
$ws.Range("E21").Value = 12
$ws.Range("G21").Value = 1
